$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.635.07"
$ws.Range("E2").Value = "  +5.75%  "
$ws.Range("D3").Value = "'3.484.90"
$ws.Range("E3").Value = "  +6.84%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'583.84"
$ws.Range("E5").Value = "  +6.75%  "
$ws.Range("D6").Value = "'158.65"
$ws.Range("E6").Value = "  +7.08%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'3.488.05"
$ws.Range("E8").Value = "  +6.59%  "
$ws.Range("D9").Value = "'0.535"
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("E11").Value = "  +7.33%  "
$ws.Range("D12").Value = "'0.442"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").Value = "'4.076.40"
$ws.Range("E13").Value = "  +6.69%  "
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "'0.0000190"
$ws.Range("E15").Value = "  +7.52%  "
$ws.Range("D16").Value = "'27.75"
$ws.Range("E16").Value = "  +4.75%  "
$ws.Range("D17").Value = "'64.635.29"
$ws.Range("E17").Value = "  +5.85%  "
$ws.Range("D18").Value = "'3.451.17"
$ws.Range("E18").Value = "  +6.05%  "
$ws.Range("D19").Value = "'6.46"
$ws.Range("E19").Value = "  +2.27%  "
$ws.Range("D20").Value = "'14.39"
$ws.Range("E20").Value = "  +6.86%  "
$ws.Range("D21").Value = "'398.73"
$ws.Range("E21").Value = "  +4.91%  "
$ws.Range("D22").Value = "'8.57"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").Value = "'0.547"
$ws.Range("E23").Value = "  +2.61%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").Value = "'72.19"
$ws.Range("E25").Value = "  +2.87%  "
$ws.Range("D26").Value = "'0.0000111"
$ws.Range("E26").Value = "  +19.35%  "
$ws.Range("D27").Value = "'9.56"
$ws.Range("E27").Value = "  +10.59%  "
$ws.Range("E28").Value = "  +6.09%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("E30").Value = "  +12.43%  "
$ws.Range("D31").Value = "'6.72"
$ws.Range("E31").Value = "  +7.94%  "
$ws.Range("E32").Value = "  +6.04%  "
$ws.Range("D33").Value = "'5.86"
$ws.Range("E33").Value = "  +7.95%  "
$ws.Range("D34").Value = "'23.93"
$ws.Range("E34").Value = "  +5.64%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").Value = "'6.97"
$ws.Range("E36").Value = "  +4.52%  "
$ws.Range("E37").Value = "  +4.85%  "
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("D39").Value = "'28.57"
$ws.Range("E39").Value = "  +8.05%  "
$ws.Range("D40").Value = "'0.0786"
$ws.Range("E40").Value = "  +8.80%  "
$ws.Range("D41").Value = "'1.90"
$ws.Range("E41").Value = "  +9.50%  "
$ws.Range("D42").Value = "'2.899.66"
$ws.Range("E42").Value = "  +3.11%  "
$ws.Range("D43").Value = "'0.0325"
$ws.Range("E43").Value = "  +3.78%  "
$ws.Range("D44").Value = "'0.788"
$ws.Range("E44").Value = "  +7.18%  "
$ws.Range("D45").Value = "'4.44"
$ws.Range("E45").Value = "  +3.56%  "
$ws.Range("D46").Value = "'42.24"
$ws.Range("E46").Value = "  +5.10%  "
$ws.Range("D47").Value = "'1.12"
$ws.Range("E47").Value = "  +10.11%  "
$ws.Range("D48").Value = "'3.528.43"
$ws.Range("E48").Value = "  +6.84%  "
$ws.Range("D49").Value = "'22.84"
$ws.Range("E49").Value = "  +5.59%  "
$ws.Range("D50").Value = "'2.16"
$ws.Range("E50").Value = "  +23.58%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "'299.98"
$ws.Range("E51").Value = "  +7.67%  "
